$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 13 ("Docentes responsáveis:" value row, holding the professor
# name under the wrong label) is removed entirely; everything below shifts
# up by one row. Deleting the row preserves each remaining row's own
# height/format, which already lines up with the target layout for most
# rows.
$ws.Rows.Item(13).Delete() | Out-Null

# After the shift, a handful of B/C text values need to be corrected so the
# right body text sits under the right heading.
$ws.Range("B10").Value = "2143261 - André Luis Ferraz"
$ws.Range("C10").Value = "2143261 - André Luis Ferraz"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Value = "01/01/2018"
$ws.Range("C15").Value = "01/01/2018"

$ws.Range("B18").Value = "2143261 - André Luis Ferraz"
$ws.Range("C18").Value = "2143261 - André Luis Ferraz"

$ws.Range("B19").Value = "A avaliação será feita por meio de provas escritas."
$ws.Range("C19").Value = "A avaliação será feita por meio de provas escritas."

$ws.Range("B20").Value = "A nota final (NF) será calculada da seguintes maneira: NF=(P1+P2)/2"
$ws.Range("C20").Value = "A nota final (NF) será calculada da seguintes maneira: NF=(P1+P2)/2"

$ws.Range("B21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) será calculada como MR=(NF=PR)/2"
$ws.Range("C21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) será calculada como MR=(NF=PR)/2"
